# Auto-generated edit script
# Updates numeric profit/price figures across the Jenova_Profits leve tables
# (currentAveragePrice / LevePrice / LeveProfit columns) per the scheduled
# market-data refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1265.16
$ws.Range("J17").Value = 1493.95
$ws.Range("L17").Value = 4481.85
$ws.Range("N17").Value = -4817.85
$ws.Range("H40").Value = 7707.923
$ws.Range("I40").Value = 6201
$ws.Range("J40").Value = 8649.75
$ws.Range("K40").Value = 6201
$ws.Range("L40").Value = 8649.75
$ws.Range("M40").Value = -6026
$ws.Range("N40").Value = -8999.75
$ws.Range("H94").Value = 2247.5557
$ws.Range("I94").Value = 1903.5
$ws.Range("K94").Value = 1903.5
$ws.Range("M94").Value = -1452.5
$ws.Range("H134").Value = 69999.89999999999
$ws.Range("J134").Value = 69999.89999999999
$ws.Range("L134").Value = 69999.89999999999
$ws.Range("N134").Value = -80139.89999999999
$ws.Range("H137").Value = 366258.84
$ws.Range("I137").Value = 271965.8
$ws.Range("K137").Value = 815897.3999999999
$ws.Range("M137").Value = -813347.3999999999
$ws.Range("H138").Value = 5773.841
$ws.Range("J138").Value = 5775.1143
$ws.Range("L138").Value = 17325.3429
$ws.Range("N138").Value = -27605.3429
$ws.Range("H141").Value = 2857.359
$ws.Range("I141").Value = 1481.2333
$ws.Range("K141").Value = 4443.699900000001
$ws.Range("M141").Value = 736.3000999999995
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6336.875
$ws.Range("I32").Value = 4746
$ws.Range("K32").Value = 4746
$ws.Range("M32").Value = -4459
$ws.Range("H45").Value = 3213.9546
$ws.Range("I45").Value = 3035.8667
$ws.Range("K45").Value = 3035.8667
$ws.Range("M45").Value = -2658.8667
$ws.Range("H74").Value = 1292.1923
$ws.Range("I74").Value = 1279.125
$ws.Range("J74").Value = 1449
$ws.Range("K74").Value = 1279.125
$ws.Range("L74").Value = 1449
$ws.Range("M74").Value = -405.125
$ws.Range("N74").Value = -3197
$ws.Range("H77").Value = 1292.1923
$ws.Range("I77").Value = 1279.125
$ws.Range("J77").Value = 1449
$ws.Range("K77").Value = 6395.625
$ws.Range("L77").Value = 7245
$ws.Range("M77").Value = -2027.625
$ws.Range("N77").Value = -15981
$ws.Range("H97").Value = 4800.846
$ws.Range("I97").Value = 5400.3
$ws.Range("J97").Value = 2802.6667
$ws.Range("K97").Value = 5400.3
$ws.Range("L97").Value = 2802.6667
$ws.Range("M97").Value = -4904.3
$ws.Range("N97").Value = -3794.6667
$ws.Range("H132").Value = 3327.4358
$ws.Range("I132").Value = 2070.1614
$ws.Range("K132").Value = 6210.4842
$ws.Range("M132").Value = -3680.4842
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2545.8262
$ws.Range("I94").Value = 1121.1904
$ws.Range("J94").Value = 17504.5
$ws.Range("K94").Value = 1121.1904
$ws.Range("L94").Value = 17504.5
$ws.Range("M94").Value = -670.1904
$ws.Range("N94").Value = -18406.5
$ws.Range("H134").Value = 36337.305
$ws.Range("I134").Value = 5041.769
$ws.Range("J134").Value = 56679.4
$ws.Range("K134").Value = 15125.307
$ws.Range("L134").Value = 170038.2
$ws.Range("M134").Value = -12590.307
$ws.Range("N134").Value = -175108.2
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 306058.78
$ws.Range("I31").Value = 335227.06
$ws.Range("K31").Value = 335227.06
$ws.Range("M31").Value = -334932.06
$ws.Range("H34").Value = 306058.78
$ws.Range("I34").Value = 335227.06
$ws.Range("K34").Value = 335227.06
$ws.Range("M34").Value = -335025.06
$ws.Range("H58").Value = 248807.7
$ws.Range("I58").Value = 591783.2
$ws.Range("K58").Value = 591783.2
$ws.Range("M58").Value = -591580.2
$ws.Range("H132").Value = 4949.3213
$ws.Range("I132").Value = 4007.5264
$ws.Range("J132").Value = 6937.5557
$ws.Range("K132").Value = 12022.5792
$ws.Range("L132").Value = 20812.6671
$ws.Range("M132").Value = -9492.5792
$ws.Range("N132").Value = -25872.6671
$ws.Range("H134").Value = 412434.6
$ws.Range("I134").Value = 325439.66
$ws.Range("K134").Value = 976318.98
$ws.Range("M134").Value = -973783.98
$ws.Range("H135").Value = 64476.844
$ws.Range("J135").Value = 64476.844
$ws.Range("L135").Value = 64476.844
$ws.Range("N135").Value = -74616.844
$ws.Range("H136").Value = 248807.7
$ws.Range("I136").Value = 591783.2
$ws.Range("K136").Value = 1775349.6
$ws.Range("M136").Value = -1772799.6
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 715017.5
$ws.Range("J92").Value = 1144
$ws.Range("L92").Value = 3432
$ws.Range("N92").Value = -5928
$ws.Range("H137").Value = 2627.7778
$ws.Range("I137").Value = 3225.5
$ws.Range("J137").Value = 1432.3334
$ws.Range("K137").Value = 9676.5
$ws.Range("L137").Value = 4297.0002
$ws.Range("M137").Value = -4576.5
$ws.Range("N137").Value = -14497.0002
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 4999
$ws.Range("I41").Value = 4999
$ws.Range("K41").Value = 4999
$ws.Range("M41").Value = -4644
$ws.Range("H70").Value = 7990.8184
$ws.Range("I70").Value = 7842.857
$ws.Range("J70").Value = 8249.75
$ws.Range("K70").Value = 7842.857
$ws.Range("L70").Value = 8249.75
$ws.Range("M70").Value = -7572.857
$ws.Range("N70").Value = -8789.75
$ws.Range("H73").Value = 7990.8184
$ws.Range("I73").Value = 7842.857
$ws.Range("J73").Value = 8249.75
$ws.Range("K73").Value = 7842.857
$ws.Range("L73").Value = 8249.75
$ws.Range("M73").Value = -6906.857
$ws.Range("N73").Value = -10121.75
$ws.Range("H102").Value = 2319.4666
$ws.Range("I102").Value = 1104.9375
$ws.Range("K102").Value = 1104.9375
$ws.Range("M102").Value = 517.0625
$ws.Range("H106").Value = 11333.333
$ws.Range("J106").Value = 11333.333
$ws.Range("L106").Value = 11333.333
$ws.Range("N106").Value = -13857.333
$ws.Range("H126").Value = 7332.8
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 7332.8
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 21998.4
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -26938.4
$ws.Range("H132").Value = 728361.4
$ws.Range("I132").Value = 913248.5600000001
$ws.Range("J132").Value = 147287.42
$ws.Range("K132").Value = 2739745.68
$ws.Range("L132").Value = 441862.26
$ws.Range("M132").Value = -2737215.68
$ws.Range("N132").Value = -446922.26
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1517.7273
$ws.Range("I16").Value = 1436.2
$ws.Range("K16").Value = 1436.2
$ws.Range("M16").Value = -1266.2
$ws.Range("H82").Value = 2924.1428
$ws.Range("J82").Value = 1987.6666
$ws.Range("L82").Value = 1987.6666
$ws.Range("N82").Value = -2709.6666
$ws.Range("H85").Value = 2924.1428
$ws.Range("J85").Value = 1987.6666
$ws.Range("L85").Value = 1987.6666
$ws.Range("N85").Value = -4483.6666
$ws.Range("H93").Value = 2286.1538
$ws.Range("I93").Value = 2152.7896
$ws.Range("J93").Value = 2648.1428
$ws.Range("K93").Value = 2152.7896
$ws.Range("L93").Value = 2648.1428
$ws.Range("M93").Value = -904.7896000000001
$ws.Range("N93").Value = -5144.1428
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 36500
$ws.Range("J75").Value = 36500
$ws.Range("L75").Value = 36500
$ws.Range("N75").Value = -38372
$ws.Range("H78").Value = 36500
$ws.Range("J78").Value = 36500
$ws.Range("L78").Value = 109500
$ws.Range("N78").Value = -118860
$ws.Range("H86").Value = 66333.336
$ws.Range("J86").Value = 66333.336
$ws.Range("L86").Value = 66333.336
$ws.Range("N86").Value = -68579.336
$ws.Range("H89").Value = 66333.336
$ws.Range("J89").Value = 66333.336
$ws.Range("L89").Value = 331666.68
$ws.Range("N89").Value = -342898.68
$ws.Range("H96").Value = 201379.8
$ws.Range("I96").Value = 333966.34
$ws.Range("K96").Value = 333966.34
$ws.Range("M96").Value = -332593.34
$ws.Range("H100").Value = 905.86365
$ws.Range("I100").Value = 704.5
$ws.Range("K100").Value = 1409
$ws.Range("M100").Value = -868
$ws.Range("H122").Value = 34486636
$ws.Range("I122").Value = 52634130
$ws.Range("K122").Value = 157902390
$ws.Range("M122").Value = -157899940
$ws.Range("H132").Value = 24036.766
$ws.Range("I132").Value = 1999.5143
$ws.Range("J132").Value = 88312.086
$ws.Range("K132").Value = 5998.5429
$ws.Range("L132").Value = 264936.258
$ws.Range("M132").Value = -3468.5429
$ws.Range("N132").Value = -269996.258
$ws.Range("H136").Value = 890378.5600000001
$ws.Range("I136").Value = 1159488.9
$ws.Range("K136").Value = 3478466.7
$ws.Range("M136").Value = -3475916.7
